## Studenci.xlsx - "Add files via upload" edit
## 1) Wynagrodzenie sheet: fill in missing 2018 salary (C25) and the two
##    derived year-over-year inflation formulas (D24, D25) that become
##    computable once C25 is known.
## 2) Inflacja sheet: collapse the monthly CPI series (2012-1 .. 2021-12)
##    into a yearly summary (2012 .. 2021) in columns A:B, and leave a
##    blank helper column C next to the new yearly rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Wynagrodzenie
# ---------------------------------------------------------------------
$wyn = $wb.Worksheets.Item("Wynagrodzenie")

$wyn.Range("C25").Value = 4585
$wyn.Range("D24").Formula = "=(C24-C25)/C25*100"
$wyn.Range("D25").Formula = "=(C25-C26)/C26*100"

$wyn.Range("F26").Select()

# ---------------------------------------------------------------------
# 2) Inflacja
# ---------------------------------------------------------------------
$infl = $wb.Worksheets.Item("Inflacja")

$years = @(2021, 2020, 2019, 2018, 2017, 2016, 2015, 2014, 2013, 2012)
$values = @(8.5999999999999943, 2.4000000000000057, 3.4000000000000057, 1.0999999999999943, 2.0999999999999943, 0.79999999999999716, -0.5, -1, 0.70000000000000284, 2.4000000000000057)

# Clear out the old monthly rows (A2:B121) first.
$infl.Range("A2:C121").ClearContents()

for ($i = 0; $i -lt $years.Length; $i++) {
    $r = $i + 2
    $infl.Cells.Item($r, 1).Value = $years[$i]
    $infl.Cells.Item($r, 2).Value = $values[$i]
}

$infl.Range("C2:C11").NumberFormat = "0.0"
$infl.Range("C2:C11").WrapText = $true
$infl.Range("C2:C11").VerticalAlignment = -4108

$infl.Range("B12:B121").NumberFormat = "0.0"

$infl.Range("G8").Select()
